$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Group number"
$ws.Range("B1").Value = "Member"
$ws.Range("C1").Value = "1. Conceptual Design"
$ws.Range("D1").Value = "2. Logical Design"
$ws.Range("E1").Value = "3. Implementation"
$ws.Range("F1").Value = "4. Database Instance"
$ws.Range("G1").Value = "5. SQL Table Modifications"
$ws.Range("H1").Value = "6. SQL Data Queries"
$ws.Range("I1").Value = "7. SQL Programming"
$ws.Range("J1").Value = "8. Java Database Access"
$ws.Range("K1").Value = "AVERAGE"

# Member names
$ws.Range("B2").Value = "Bjarne Larsen"
$ws.Range("B3").Value = "Marcus Lemser"
$ws.Range("B4").Value = "Maximillian Mortesen"
$ws.Range("B5").Value = "Oscar Bjerregaard"
$ws.Range("B6").Value = "Tobias Frederiksen"

# Autofit columns A:K to match bestFit widths from the diff
$ws.Range("A1:K6").Columns.AutoFit()

# Update selection to match target state
$ws.Range("C10").Select()
